$wb = $excel.ActiveWorkbook

# --- Add the three new worksheets in the correct tab order ---
$classesSheet = $wb.Worksheets.Item("Classes")
$todoList = $wb.Worksheets.Add($null, $classesSheet)
$todoList.Name = "TodoList"

$props = $wb.Worksheets.Add($null, $todoList)
$props.Name = "Properties"

$genOrder = $wb.Worksheets.Item("Generation Order")
$sheet3 = $wb.Worksheets.Add($null, $genOrder)
$sheet3.Name = "Sheet3"

# --- Populate "TodoList" sheet ---
$todoList.Range('A1').Value = 'Features to Code'
$todoList.Range('E1').Value = 'Features finished'
$todoList.Range('N1').Value = 'Features:'
$todoList.Range('N2').Value = '"Generates random maps with custom biomes"'
$todoList.Range('A3').Value = 'Generation:'
$todoList.Range('O3').Value = '"Biomes will influence others nearby, sharing properties and transitioning smoothly"'
$todoList.Range('A4').Value = 'Add "do not start" biomes '
$todoList.Range('N4').Value = '"Adds terrain features such as towns, dungeons, lairs, roaming monsters, roads and rivers"'
$todoList.Range('A5').Value = 'Add River Generation'
$todoList.Range('O5').Value = '"Gives information of surrounding hexes based on what the players could see"'
$todoList.Range('A6').Value = 'Add town and landmark Gen'
$todoList.Range('N6').Value = '"Creates a auto-generated d100 encounter list for each hex"'
$todoList.Range('A7').Value = 'Add Road Generation'
$todoList.Range('O7').Value = '"Create custom encounter lists for worldobjects that the generator will insert"'
$todoList.Range('A8').Value = 'Add River and Road nodes'
$todoList.Range('A9').Value = 'Add Monster Generation'
$todoList.Range('A10').Value = 'Add biome modificators from file'
$todoList.Range('A11').Value = 'Add biome modificators from worldobjects'
$todoList.Range('A12').Value = 'Add Custom Encounters to replace standard encounters'
$todoList.Range('N12').Value = 'FUNCTIONALITY'
$todoList.Range('N13').Value = 'Rivers'
$todoList.Range('Q13').Value = 'Roads'
$todoList.Range('N14').Value = 'Link different hexes together'
$todoList.Range('Q14').Value = 'Link towns together'
$todoList.Range('A15').Value = 'Editing:'
$todoList.Range('N15').Value = 'Act as boundaries to monsters'
$todoList.Range('Q15').Value = 'Generate signposts for directions'
$todoList.Range('A16').Value = 'Save/Load hexmap with serialization'
$todoList.Range('N16').Value = 'Increase travel cost'
$todoList.Range('Q16').Value = 'Generate crossings on rivers (remove travel cost)'
$todoList.Range('A22').Value = 'Outputting:'

# --- Populate "Properties" sheet ---
$props.Range('E1').Value = 'Chance is RandInt(1,100), if result < chance, thing happens.'
$props.Range('N1').Value = 'Visibility affects how well it is seen within and nearby its home hex. Large visibility means it''ll appear on more encounter tables as "evidence"'
$props.Range('O2').Value = 'EVIDENCE'
$props.Range('R2').Value = '(Lair)'
$props.Range('S2').Value = 'no lairs, create 0 nomadicchance monsters to always generate a lair.'
$props.Range('S3').Value = 'if additional flavor required, use dungeons instead.'
$props.Range('B4').Value = 'Monster:'
$props.Range('P4').Value = 'Town:'
$props.Range('C5').Value = '<Name>,<O_DMNotes>,<O_MaxAmount >=1>, <RoamChance>-<O_RoamRadius>,<BiomeName1..X>-<%>…'
$props.Range('Q5').Value = 'Name'
$props.Range('Q6').Value = 'Description'
$props.Range('C7').Value = 'Name:'
$props.Range('Q7').Value = 'SpawnChance: <biomename>-X'
$props.Range('C8').Value = 'Description:'
$props.Range('Q8').Value = 'Visibility:'
$props.Range('C9').Value = 'NomadChance:'
$props.Range('E9').Value = 'Lair or nomadic'
$props.Range('Q9').Value = 'Connectivity:'
$props.Range('S9').Value = 'The resource it uses to create a road to nearby town(s)'
$props.Range('C10').Value = 'RoamingRadius:'
$props.Range('E10').Value = 'applies to both nomadic and lairs'
$props.Range('I10').Value = 'Affected by travel cost, does a "flood" search to find radius'
$props.Range('C11').Value = 'SpawnChance: <biomename>-X, <b2>-X'
$props.Range('H11').Value = 'X% chance PER hex of that biome type.'
$props.Range('C12').Value = 'Visibility:'
$props.Range('D12').Value = 'Likelihood of encountering tracks or evidence, not the monster itself.'
$props.Range('C13').Value = 'RoamType:'
$props.Range('D13').Value = 'Fly or Other'
$props.Range('F13').Value = 'cannot ''swim'' down rivers, as they are not biomes.'
$props.Range('F14').Value = 'Travelcost affects ''land'' creatures'
$props.Range('C15').Value = 'EncounterChance:'
$props.Range('E15').Value = 'X'
$props.Range('B17').Value = 'Biome'
$props.Range('P17').Value = 'Dungeon'
$props.Range('C18').Value = 'Name:'
$props.Range('Q18').Value = 'Name'
$props.Range('C19').Value = 'Description'
$props.Range('E19').Value = 'Use full sentences, 3 maximum probably.'
$props.Range('Q19').Value = 'Description'
$props.Range('S19').Value = 'appends any monsters as ''patrols'''
$props.Range('C20').Value = 'Color'
$props.Range('Q20').Value = 'SpawnChance:'
$props.Range('S20').Value = '<biomename>-X'
$props.Range('C21').Value = 'height'
$props.Range('Q21').Value = 'MonsterChance: '
$props.Range('S21').Value = '<monstername>-X, <monstername>-X'
$props.Range('W21').Value = 'lair monsters that can be present'
$props.Range('C22').Value = 'travelcost'
$props.Range('Q22').Value = 'Visibility'
$props.Range('W22').Value = 'how dominant evidence of these tracks are'
$props.Range('C23').Value = 'spotdistance'
$props.Range('E23').Value = 'encounters will be encountered at this distance within the hex (ft)'
$props.Range('C24').Value = 'riverorigin'
$props.Range('F24').Value = 'aka, stealth rolls should be rolled at this distance, failure = this distance spotted.'
$props.Range('C25').Value = 'riverend'
$props.Range('F25').Value = 'spot should be a diceroll string'
$props.Range('B28').Value = 'BiomeModifier:'
$props.Range('D28').Value = 'ADDS to modified biome'
$props.Range('P28').Value = 'Landmark'
$props.Range('C29').Value = 'Name: '
$props.Range('D29').Value = '(Prefix)'
$props.Range('Q29').Value = 'Name'
$props.Range('C30').Value = 'Description:'
$props.Range('D30').Value = 'Adds, use 1 sentence./'
$props.Range('G30').Value = 'Some minor modifiers such as river origins only add description, not name.'
$props.Range('Q30').Value = 'Description'
$props.Range('C31').Value = 'Color: '
$props.Range('D31').Value = 'adds this component to each'
$props.Range('Q31').Value = 'SpawnChance:'
$props.Range('C32').Value = 'height: '
$props.Range('D32').Value = 'adds if positive, takes if negative'
$props.Range('Q32').Value = 'Visibility'
$props.Range('C33').Value = 'travelcost'
$props.Range('D33').Value = 'as above'
$props.Range('C34').Value = 'spotdistance '
$props.Range('D34').Value = 'as above'
$props.Range('G34').Value = 'flat integer'
$props.Range('C35').Value = 'validbiomes:'
$props.Range('G35').Value = 'will only apply this biome to valid biomes '
$props.Range('C36').Value = 'origin: '
$props.Range('D36').Value = '<biomename>-X'
$props.Range('F36').Value = 'biome names csv'
$props.Range('H36').Value = 'searches via name, regardless of type.'
$props.Range('C37').Value = 'riverorigin'
$props.Range('E37').Value = 'affects river generation'
$props.Range('I37').Value = 'Biomes only'
$props.Range('C38').Value = 'riverend'
$props.Range('E38').Value = 'affects river generation'
$props.Range('I38').Value = 'PRE-SET MODIFIERS THAT OCCUR WITHIN PROGRAM'
$props.Range('P38').Value = 'Region'
$props.Range('T38').Value = 'After generation, finds clusters of valid biomes and labels them with a region'
$props.Range('Q39').Value = 'Name'
$props.Range('U39').Value = 'This is mainly for lore'
$props.Range('Q40').Value = 'Description'
$props.Range('Q41').Value = 'MinMax'
$props.Range('R41').Value = '<Min,Max>'
$props.Range('Q42').Value = 'MinSize'
$props.Range('R42').Value = 'X'
$props.Range('T42').Value = 'How many hexes of validbiomes together before this is considered a valid region location.'
$props.Range('Q43').Value = 'validbiomes:'

# --- Populate "Sheet3" sheet ---
$sheet3.Range('B1').Value = 'Road Generator'
$sheet3.Range('B2').Value = 'Sort Towns by largest connectivity TownQ'
$sheet3.Range('B3').Value = 'While initialpath=false'
$sheet3.Range('C4').Value = 'Get largest Town (or next largest)'
$sheet3.Range('C5').Value = 'If can reach nearest town'
$sheet3.Range('H5').Value = 'Dijkstra''s'
$sheet3.Range('D6').Value = 'ROADPath to nearest town '
$sheet3.Range('D7').Value = 'InitialPath=true'
$sheet3.Range('D8').Value = '(do not subtract connectivity for initial path)'
$sheet3.Range('B10').Value = 'While contains towns in townQ'
$sheet3.Range('C11').Value = 'Get next largest connectivity'
$sheet3.Range('C12').Value = 'If can reach nearest town or roadnode'
$sheet3.Range('H12').Value = 'Dijkstra''s'
$sheet3.Range('I12').Value = 'search through connections for hexes containing towns, if they do'
$sheet3.Range('D13').Value = 'ROADconnect to nearest town OR roadnode'
$sheet3.Range('D14').Value = 'subtract connectivity from BOTH'
$sheet3.Range('D15').Value = 'update townQ'
$sheet3.Range('C16').Value = 'else'
$sheet3.Range('D17').Value = 'remove town from townQ'

# --- Add new rows to "Generation Order" sheet ---
$genOrder.Range('D7').Value = 'Valid Starts:'
$genOrder.Range('E8').Value = 'Elevated Water sources'
$genOrder.Range('E9').Value = 'Mountains'
$genOrder.Range('D11').Value = 'Valid Ends:'
$genOrder.Range('E12').Value = 'Low water sources'

# --- Add new row 20 to "Classes" sheet ---
$classesSheet.Range('C20').Value = 'Encounters '
$classesSheet.Range('E20').Value = '0..* '
$classesSheet.Range('F20').Value = 'uses this list when adding encounters to the d100 list, replaces standard "there are ghouls here"'

# --- Update selections on each sheet ---
$importSheet = $wb.Worksheets.Item("Import File Formats")
$importSheet.Range("C18").Select()

$classesSheet.Range("F21").Select()

$todoList.Range("A11").Select()

$props.Range("T43").Select()

$genOrder.Range("E13").Select()

$sheet3.Range("L27").Select()

